# "NIH Terminations Data Analysis Final 5.18.25"
#
# This workbook's single data table (terminated_grants!A1:E12) gets re-sorted
# descending by column D ("Total Including Terminated") instead of column E
# ("Percent Termination"), and a new blank worksheet named "Sheet1" is added
# after the existing "terminated_grants" sheet. The active-cell selection on
# terminated_grants also moves from C10 to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the new, empty "Sheet1" worksheet right after terminated_grants ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "Sheet1"

# Keep terminated_grants as the active/tab-selected sheet
$ws.Activate()

# --- Re-sort the data table by Total Including Terminated (col D), desc ---
$sortRange = $ws.Range("A1:E13")
$sortKey = $ws.Range("D1:D13")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 2, 0, 0) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Apply()

# The sort leaves the cell that used to hold the (blank) "C" value for the
# "G" grant family without the highlighted (yellow) formatting that the rest
# of that row carries after moving to row 12 - restore it to match its
# siblings in the row.
$ws.Range("C12").Interior.Color = $ws.Range("B12").Interior.Color

# --- Update the saved selection/active cell on terminated_grants ---
$ws.Range("D7").Select() | Out-Null
